$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 104: correct the date/time serial value in column A ---
$ws.Range("A104").Value = 45485.2916666667

# --- Row 105: new row appended (results from the R script) ---

# A105 reuses the same date/time style as A104 -- copy the formatted cell
# down first (so the datetime number format + font come along), then set
# the new serial value.
$ws.Range("A104").Copy($ws.Range("A105"))
$ws.Range("A105").Value = 45488.6178819444

$ws.Range("B105").Value = 8400
$ws.Range("C105").Value = 6.05999994277954
$ws.Range("D105").Value = 6
$ws.Range("E105").Value = 6.05999994277954
$ws.Range("F105").Value = 6

# G105/H105 hold plain numeric-/label-looking text in this sheet (shared
# strings, t="s"), not numbers. Assigning the literal string directly would
# get auto-coerced to a number by Excel, so instead build the text via a
# formula that evaluates to a string, then collapse the formula down to its
# literal (text) value in place with a values-only paste. This keeps the
# cell's number format/style untouched (still the default style).
$ws.Range("G105").Formula = '="6"'
$ws.Range("G105").Copy()
$ws.Range("G105").PasteSpecial(-4163)

$ws.Range("H105").Formula = '="PAL.MI"'
$ws.Range("H105").Copy()
$ws.Range("H105").PasteSpecial(-4163)
